# Add files via upload
# The sheet originally held a 4-column product table in A1:D4 (ID, Name, Price, Stock).
# The new upload adds a pandas-style integer index column in front, so the
# original columns shift from A:D to B:E, and the data itself is replaced by
# a new product list (Apple/Banana/Cherry with whole-number prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole existing table one column to the right to make room for the
# new leading index column (this carries the header style along with it).
$ws.Range("A1:A4").Insert(-4161)

# New index column (pandas RangeIndex 0,1,2) using the same style the header
# row already uses (border + bold/center) per the style id seen in the diff.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("B1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# Replace the product data with the new rows.
$ws.Range("B2").Value = 101
$ws.Range("C2").Value = "Apple"
$ws.Range("D2").Value = 40
$ws.Range("E2").Value = 100

$ws.Range("B3").Value = 102
$ws.Range("C3").Value = "Banana"
$ws.Range("D3").Value = 35
$ws.Range("E3").Value = 75

$ws.Range("B4").Value = 103
$ws.Range("C4").Value = "Cherry"
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 50
